$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '45.405.06'
$ws.Range("E2").Value = '  +6.18%  '

# Row 3
$ws.Range("D3").Value = '2.365.57'
$ws.Range("E3").Value = '  +2.21%  '

# Row 4
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").Value = '109.48'
$ws.Range("E5").Value = '  +2.19%  '

# Row 6
$ws.Range("D6").Value = '309.76'
$ws.Range("E6").Value = '  -0.66%  '

# Row 7
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("E8").Value = '  -0.18%  '

# Row 9
$ws.Range("D9").Value = '0.618'
$ws.Range("E9").Value = '  +1.84%  '

# Row 10
$ws.Range("D10").Value = '41.24'
$ws.Range("E10").Value = '  +2.53%  '

# Row 11
$ws.Range("D11").Value = '0.0919'
$ws.Range("E11").Value = '  +0.46%  '

# Row 12
$ws.Range("D12").Value = '8.49'
$ws.Range("E12").Value = '  +1.41%  '

# Row 13
$ws.Range("E13").Value = '  +2.14%  '

# Row 14
$ws.Range("D14").Value = '0.986'
$ws.Range("E14").Value = '  -0.57%  '

# Row 15
$ws.Range("D15").Value = '2.724.54'
$ws.Range("E15").Value = '  +2.31%  '

# Row 16
$ws.Range("D16").Value = '15.39'
$ws.Range("E16").Value = '  +0.19%  '

# Row 17
$ws.Range("D17").Value = '2.390.28'
$ws.Range("E17").Value = '  +3.57%  '

# Row 18
$ws.Range("D18").Value = '45.365.80'
$ws.Range("E18").Value = '  +5.55%  '

# Row 19
$ws.Range("D19").Value = '7.32'
$ws.Range("E19").Value = '  -2.18%  '

# Row 20
$ws.Range("E20").Value = '  +0.58%  '

# Row 21
$ws.Range("D21").Value = '13.69'
$ws.Range("E21").Value = '  +4.76%  '

# Row 22
$ws.Range("D22").Value = '73.39'
$ws.Range("E22").Value = '  -0.33%  '

# Row 23
$ws.Range("E23").Value = '  -0.42%  '

# Row 24
$ws.Range("D24").Value = '258.96'
$ws.Range("E24").Value = '  -2.64%  '

# Row 25
$ws.Range("E25").Value = '  +3.06%  '

# Row 26
$ws.Range("E26").Value = '  -0.53%  '

# Row 27
$ws.Range("D27").Value = '11.17'
$ws.Range("E27").Value = '  +1.30%  '

# Row 28
$ws.Range("D28").Value = '7.36'
$ws.Range("E28").Value = '  -5.12%  '

# Row 29
$ws.Range("D29").Value = '2.36'
$ws.Range("E29").Value = '  +3.06%  '

# Row 30
$ws.Range("D30").Value = '0.0973'
$ws.Range("E30").Value = '  +11.11%  '

# Row 31
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '38.04'
$ws.Range("E31").Value = '  -1.28%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '22.37'
$ws.Range("E32").Value = '  -0.34%  '

# Row 33
$ws.Range("D33").Value = '169.36'
$ws.Range("E33").Value = '  +1.81%  '

# Row 34
$ws.Range("E34").Value = '  +6.38%  '

# Row 35
$ws.Range("E35").Value = '  +0.36%  '

# Row 36
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.117'
$ws.Range("E36").Value = '  +4.01%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.84'
$ws.Range("E37").Value = '  +3.31%  '

# Row 38
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '2.97'
$ws.Range("E38").Value = '  +4.38%  '

# Row 39
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = '3.94'
$ws.Range("E39").Value = '  +6.96%  '

# Row 40
$ws.Range("D40").Value = '0.0357'
$ws.Range("E40").Value = '  -0.18%  '

# Row 41
$ws.Range("E41").Value = '  +8.48%  '

# Row 42
$ws.Range("D42").Value = '99.22'
$ws.Range("E42").Value = '  -4.48%  '

# Row 43
$ws.Range("D43").Value = '0.233'
$ws.Range("E43").Value = '  -0.24%  '

# Row 44
$ws.Range("D44").Value = '69.91'
$ws.Range("E44").Value = '  -1.61%  '

# Row 45
$ws.Range("D45").Value = '13.03'
$ws.Range("E45").Value = '  +1.08%  '

# Row 46
$ws.Range("E46").Value = '  -0.26%  '

# Row 47
$ws.Range("D47").Value = '83.16'
$ws.Range("E47").Value = '  +8.04%  '

# Row 48
$ws.Range("D48").Value = '112.54'
$ws.Range("E48").Value = '  -0.08%  '

# Row 49
$ws.Range("D49").Value = '5.49'
$ws.Range("E49").Value = '  +4.43%  '

# Row 50
$ws.Range("D50").Value = '9.15'
$ws.Range("E50").Value = '  +3.93%  '

# Row 51
$ws.Range("D51").Value = '1.674.88'
$ws.Range("E51").Value = '  +1.31%  '

